$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column C, matching the style of the existing header cells
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "Electrode Locations"

# Re-sort the electrode-location data (rows 2-42) by electrode location
# (letter prefix, then numeric suffix) and populate the new Electrode
# Locations column at the same time.
$data = @(
    ,@("A3_bipolar_10V_5Hz.txt", 91.820615, "A3")
    ,@("B15_bipolar_10V_5Hz.txt", 86.09393799999999, "B15")
    ,@("C1_bipolar_10V_5Hz.txt", 114.275164, "C1")
    ,@("C3_bipolar_10V_5Hz.txt", 85.452034, "C3")
    ,@("C9_bipolar_10V_5Hz.txt", 79.168232, "C9")
    ,@("C13_bipolar_10V_5Hz.txt", 76.003117, "C13")
    ,@("C15_bipolar_10V_5Hz.txt", 77.636139, "C15")
    ,@("D6_bipolar_10V_5Hz.txt", 55.179571, "D6")
    ,@("D8_bipolar_10V_5Hz.txt", 57.163823, "D8")
    ,@("E1_bipolar_10V_5Hz.txt", 90.552955, "E1")
    ,@("E4_bipolar_10V_5Hz.txt", 56.439158, "E4")
    ,@("E6_bipolar_10V_5Hz.txt", 56.638996, "E6")
    ,@("E7_bipolar_10V_5Hz.txt", 64.02291700000001, "E7")
    ,@("E10_bipolar_10V_5Hz.txt", 76.550149, "E10")
    ,@("E11_bipolar_10V_5Hz.txt", 68.970427, "E11")
    ,@("E15_bipolar_10V_5Hz.txt", 68.211445, "E15")
    ,@("G4_bipolar_10V_5Hz.txt", 58.746381, "G4")
    ,@("G5_bipolar_10V_5Hz.txt", 57.222362, "G5")
    ,@("G7_bipolar_10V_5Hz.txt", 64.68097, "G7")
    ,@("G9_bipolar_10V_5Hz.txt", 76.80045200000001, "G9")
    ,@("G12_bipolar_10V_5Hz.txt", 70.175512, "G12")
    ,@("G13_bipolar_10V_5Hz.txt", 62.147668, "G13")
    ,@("I3_bipolar_10V_5Hz.txt", 56.451269, "I3")
    ,@("I5_bipolar_10V_5Hz.txt", 58.443596, "I5")
    ,@("I11_bipolar_10V_5Hz.txt", 67.817825, "I11")
    ,@("I15_bipolar_10V_5Hz.txt", 69.521496, "I15")
    ,@("J2_bipolar_10V_5Hz.txt", 60.958731, "J2")
    ,@("J13_bipolar_10V_5Hz.txt", 66.229212, "J13")
    ,@("K1_bipolar_10V_5Hz.txt", 67.989403, "K1")
    ,@("K10_bipolar_10V_5Hz.txt", 61.455299, "K10")
    ,@("K12_bipolar_10V_5Hz.txt", 65.956705, "K12")
    ,@("K14_bipolar_10V_5Hz.txt", 64.789973, "K14")
    ,@("K15_bipolar_10V_5Hz.txt", 63.899784, "K15")
    ,@("M10_bipolar_10V_5Hz.txt", 58.909885, "M10")
    ,@("M11_bipolar_10V_5Hz.txt", 71.172684, "M11")
    ,@("M15_bipolar_10V_5Hz.txt", 60.87597, "M15")
    ,@("O2_bipolar_10V_5Hz.txt", 70.41773999999999, "O2")
    ,@("O3_bipolar_10V_5Hz.txt", 86.978071, "O3")
    ,@("O5_bipolar_10V_5Hz.txt", 92.878345, "O5")
    ,@("O8_bipolar_10V_5Hz.txt", 75.892096, "O8")
    ,@("O14_bipolar_10V_5Hz.txt", 63.336604, "O14")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $ws.Cells.Item($rowNum, 1).Value = $data[$i][0]
    $ws.Cells.Item($rowNum, 2).Value = $data[$i][1]
    $ws.Cells.Item($rowNum, 3).Value = $data[$i][2]
}

Write-Host "Electrode Locations column added and rows re-sorted."
